$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header text: remove double space before "vs."
$ws.Range("E1").Value = 'Feb 2025 vs. Jan 2025 (£)'
$ws.Range("F1").Value = 'Feb 2025 vs. Jan 2025 (%)'
$ws.Range("G1").Value = 'Mar 2025 vs. Feb 2025 (£)'
$ws.Range("H1").Value = 'Mar 2025 vs. Feb 2025 (%)'

# Add new Commentary column (copy header formatting from the last existing header cell)
$ws.Range("I1").Value = 'Commentary'
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("I1").Value = 'Commentary'

$ws.Range("I2").Value = ' Salaries: Increase of £5,000 (4.0%) vs prior period, driven by a 3 additional hires (Headcount now at 55). This equates to an average increase per employee of approximately £817 (or 2.9%), demonstrating the company''s continued focus on strategic expansion while managing salary costs effectively.'
$ws.Range("I3").Value = ' Software Licenses (Mar 2025 vs Feb 2025): Consistent expenditure at a stable level with no changes in the number of licenses. This indicates efficient utilization and no need for immediate capacity expansion or cost optimization efforts within this area.'
$ws.Range("I4").Value = ' Cloud Hosting cost increased by £1,500 (8.1%) MoM, likely due to a combination of price adjustments and increased usage. To maintain cost efficiency, it is advisable to investigate potential optimization strategies such as negotiating better pricing or optimizing resource allocation. The derived monthly growth rate KPI is 8.1%.'
$ws.Range("I5").Value = ' Travel & Entertainment: A 60% increase (£1,500) from Feb to Mar 2025 indicates a significant rise in business activities or expenses related to travel and entertainment. This surge could be due to increased meetings, events, or corporate travels. To better understand this trend, a deeper analysis of the underlying factors driving this growth is recommended, such as travel frequency, event attendance, or accommodations costs. Without specific KPIs available, it''s essential to monitor these line items closely moving forward to ensure expenses remain within budget and aligned with business objectives.'
$ws.Range("I6").Value = ' Consulting Fees decreased £500 (-5.6%) in Mar 2025 compared to Feb 2025, despite an increase in Active Projects from 5 to 6. This suggests a potential decline in the average revenue per project or reduced fees from existing clients. The Revenue per Project KPI has dropped by approximately £83 (-14.6%) from Feb 2025, indicating the need for further analysis to address this decline and optimize project profitability.'

Write-Host "done"
